$d = $word.ActiveDocument

# Find the "CONTRAST ENHANCED SPECTRAL MAMMOGRAPHY REVEALED:" heading
# paragraph, then walk backwards over the run of consecutive empty
# paragraphs that precede it and collapse that run down to a single
# blank paragraph (i.e. delete all but the first of them).
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text -replace "`r", ""
    if ($txt -like "CONTRAST ENHANCED SPECTRAL MAMMOGRAPHY REVEALED:*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -gt 0) {
    $lastEmpty = $headingIndex - 1
    $firstEmpty = $lastEmpty
    while ($firstEmpty -ge 1 -and (($d.Paragraphs.Item($firstEmpty).Range.Text -replace "`r", "") -eq "")) {
        $firstEmpty = $firstEmpty - 1
    }
    $firstEmpty = $firstEmpty + 1

    # keep the first blank paragraph of the run, delete the rest
    if ($lastEmpty -gt $firstEmpty) {
        $startPara = $d.Paragraphs.Item($firstEmpty + 1)
        $endPara = $d.Paragraphs.Item($lastEmpty)
        $rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $rangeToDelete.Delete()
    }
}

# Re-assert the heading text so Word regenerates the run and drops the
# stale <w:lastRenderedPageBreak/> rendering-cache marker left over from
# the previous pagination.
$d.Content.Find.Execute("CONTRAST ENHANCED SPECTRAL MAMMOGRAPHY REVEALED:", $true, $false, $false, $false, $false, $true, 1, $false, "CONTRAST ENHANCED SPECTRAL MAMMOGRAPHY REVEALED:", 2)
